# Add missing "X" marks for intrinsics that were evaluated but not flagged:
#   ATN()  -> row 4  (Token, Parse, Eval all supported)
#   CDBL() -> row 10 (Token, Parse, Eval all supported)
#   CINT() -> row 14 (Token, Parse, Eval all supported)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "X"
$ws.Range("C4").Value = "X"
$ws.Range("D4").Value = "X"

$ws.Range("B10").Value = "X"
$ws.Range("C10").Value = "X"
$ws.Range("D10").Value = "X"

$ws.Range("B14").Value = "X"
$ws.Range("C14").Value = "X"
$ws.Range("D14").Value = "X"

# Move the active selection to A15 (matches the post-edit saved cursor position)
$ws.Range("A15").Select()
